# Applies the commit: insert a new weekly "Piña / Caramelo / Segunda" record
# row into the Macroferia Regional de Talca sheet, shifting all subsequent
# rows down by one (A1:T475 -> A1:T476).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 440; this pushes the existing rows 440:475
# down to 441:476 and copies the formatting of the row above (439), which
# already carries the date-style used throughout this data block.
$ws.Rows.Item(440).Insert()

# Populate the newly inserted row 440 with the new record.
$ws.Cells.Item(440, 1).Value = 5
$ws.Cells.Item(440, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(440, 3).Value = "Maule"
$ws.Cells.Item(440, 4).Value = 45223
$ws.Cells.Item(440, 5).Value = 7
$ws.Cells.Item(440, 6).Value = "Fruta"
$ws.Cells.Item(440, 7).Value = 100108
$ws.Cells.Item(440, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(440, 9).Value = 100108005
$ws.Cells.Item(440, 10).Value = "Piña"
$ws.Cells.Item(440, 11).Value = "Caramelo"
$ws.Cells.Item(440, 12).Value = "Segunda"
$ws.Cells.Item(440, 13).Value = 200
$ws.Cells.Item(440, 14).Value = 21000
$ws.Cells.Item(440, 15).Value = 21000
$ws.Cells.Item(440, 16).Value = 21000
$ws.Cells.Item(440, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(440, 18).Value = "Ecuador"
$ws.Cells.Item(440, 19).Value = 1500
$ws.Cells.Item(440, 20).Value = 14
